$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (B:E) ---
$ws.Columns.Item(2).ColumnWidth = 66.49869791666667
$ws.Columns.Item(3).ColumnWidth = 102.60807291666667
$ws.Columns.Item(4).ColumnWidth = 52.721354166666664
$ws.Columns.Item(5).ColumnWidth = 36.608072916666664

# --- Update family cross-reference lists (columns B/C) ---
# Assign in an order that creates each distinct value once, in a stable sequence,
# then reuses it (string de-dup) for every other matching cell.
$ws.Range("C35").Value = '[]'
$ws.Range("C36").Value = $ws.Range("C35").Value2
$ws.Range("B2").Value = '["Apparent Matrix Density","Borehole Fluid Density","Bulk Density","Bulk Density (Array)","Bulk Density Correction","Bulk Density Hydrocarbon Corrected","Core Grain Density","Corrected ZDL Density","Density Correction","Density Count Rate","Density Porosity","Density Porosity Lime","Density Porosity Sand","Density Squared","Fluid Density","Fluid Density Contrast","Gas Density","Grain Density","Matrix Density","Mud Filtrate Density","Oil Density"]'
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("B4").Value = $ws.Range("B2").Value2
$ws.Range("B5").Value = $ws.Range("B2").Value2
$ws.Range("B6").Value = $ws.Range("B2").Value2
$ws.Range("B7").Value = $ws.Range("B2").Value2
$ws.Range("B8").Value = $ws.Range("B2").Value2
$ws.Range("B9").Value = $ws.Range("B2").Value2
$ws.Range("B11").Value = $ws.Range("B2").Value2
$ws.Range("C12").Value = $ws.Range("B2").Value2
$ws.Range("C13").Value = $ws.Range("B2").Value2
$ws.Range("C14").Value = $ws.Range("B2").Value2
$ws.Range("C15").Value = $ws.Range("B2").Value2
$ws.Range("B16").Value = $ws.Range("B2").Value2
$ws.Range("C17").Value = $ws.Range("B2").Value2
$ws.Range("B19").Value = $ws.Range("B2").Value2
$ws.Range("B20").Value = $ws.Range("B2").Value2
$ws.Range("C31").Value = $ws.Range("B2").Value2
$ws.Range("B32").Value = $ws.Range("B2").Value2
$ws.Range("B33").Value = $ws.Range("B2").Value2
$ws.Range("B37").Value = $ws.Range("B2").Value2
$ws.Range("B38").Value = $ws.Range("B2").Value2
$ws.Range("B39").Value = $ws.Range("B2").Value2
$ws.Range("B40").Value = $ws.Range("B2").Value2
$ws.Range("B43").Value = $ws.Range("B2").Value2
$ws.Range("B44").Value = $ws.Range("B2").Value2
$ws.Range("B45").Value = $ws.Range("B2").Value2
$ws.Range("B47").Value = $ws.Range("B2").Value2
$ws.Range("B48").Value = $ws.Range("B2").Value2
$ws.Range("B49").Value = $ws.Range("B2").Value2
$ws.Range("B50").Value = $ws.Range("B2").Value2
$ws.Range("B51").Value = $ws.Range("B2").Value2
$ws.Range("B52").Value = $ws.Range("B2").Value2
$ws.Range("B53").Value = $ws.Range("B2").Value2
$ws.Range("B54").Value = $ws.Range("B2").Value2
$ws.Range("C2").Value = '["Neutron Porosity","Neutron Porosity Correction","Neutron Porosity Hydrocarbon Corrected","Neutron Porosity Lime","Neutron Porosity Sand","Neutron Porosity Squared","Compensated Neutron Porosity","Thermal Neutron Porosity"]'
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("C4").Value = $ws.Range("C2").Value2
$ws.Range("C5").Value = $ws.Range("C2").Value2
$ws.Range("C6").Value = $ws.Range("C2").Value2
$ws.Range("C7").Value = $ws.Range("C2").Value2
$ws.Range("C8").Value = $ws.Range("C2").Value2
$ws.Range("C9").Value = $ws.Range("C2").Value2
$ws.Range("C11").Value = $ws.Range("C2").Value2
$ws.Range("C16").Value = $ws.Range("C2").Value2
$ws.Range("C19").Value = $ws.Range("C2").Value2
$ws.Range("C20").Value = $ws.Range("C2").Value2
$ws.Range("C21").Value = $ws.Range("C2").Value2
$ws.Range("C22").Value = $ws.Range("C2").Value2
$ws.Range("C32").Value = $ws.Range("C2").Value2
$ws.Range("C33").Value = $ws.Range("C2").Value2
$ws.Range("C34").Value = $ws.Range("C2").Value2
$ws.Range("C37").Value = $ws.Range("C2").Value2
$ws.Range("C38").Value = $ws.Range("C2").Value2
$ws.Range("C39").Value = $ws.Range("C2").Value2
$ws.Range("C40").Value = $ws.Range("C2").Value2
$ws.Range("C41").Value = $ws.Range("C2").Value2
$ws.Range("C42").Value = $ws.Range("C2").Value2
$ws.Range("C43").Value = $ws.Range("C2").Value2
$ws.Range("C44").Value = $ws.Range("C2").Value2
$ws.Range("C45").Value = $ws.Range("C2").Value2
$ws.Range("C47").Value = $ws.Range("C2").Value2
$ws.Range("C48").Value = $ws.Range("C2").Value2
$ws.Range("C49").Value = $ws.Range("C2").Value2
$ws.Range("C50").Value = $ws.Range("C2").Value2
$ws.Range("C51").Value = $ws.Range("C2").Value2
$ws.Range("C52").Value = $ws.Range("C2").Value2
$ws.Range("C53").Value = $ws.Range("C2").Value2
$ws.Range("C54").Value = $ws.Range("C2").Value2
$ws.Range("C55").Value = $ws.Range("C2").Value2
$ws.Range("C56").Value = $ws.Range("C2").Value2
$ws.Range("B12").Value = '["Acoustic","Acoustic Attenuation Rate","Acoustic Normalization Factor"]'
$ws.Range("B13").Value = $ws.Range("B12").Value2
$ws.Range("B14").Value = $ws.Range("B12").Value2
$ws.Range("B15").Value = $ws.Range("B12").Value2
$ws.Range("B17").Value = $ws.Range("B12").Value2
$ws.Range("B21").Value = $ws.Range("B12").Value2
$ws.Range("B22").Value = $ws.Range("B12").Value2
$ws.Range("B31").Value = $ws.Range("B12").Value2
$ws.Range("B34").Value = $ws.Range("B12").Value2
$ws.Range("B41").Value = $ws.Range("B12").Value2
$ws.Range("B42").Value = $ws.Range("B12").Value2
$ws.Range("B55").Value = $ws.Range("B12").Value2
$ws.Range("B56").Value = $ws.Range("B12").Value2
$ws.Range("C24").Value = '["Block Porosity","Compensated Neutron Porosity","Core Porosity","Core Porosity Under Stress (Array)","Density Porosity","Density Porosity Lime","Density Porosity Sand","Effective Porosity","Fracture Porosity","Fracture porosity cutoff - Fracture Vug Workflow","Isolated Porosity","Net Porosity","Net Porosity","Net Sand Fraction","Neutron Porosity","Neutron Porosity Correction","Neutron Porosity Hydrocarbon Corrected","Neutron Porosity Lime","Neutron Porosity Sand","Neutron Porosity Squared","Open Porosity","Open porosity cutoff - Fracture Vug Workflow","Parallel Porosity","Porosity","Porosity Pc Modeling","Porosity Unclipped","Secondary effective Porosity","Secondary Porosity","Secondary porosity cutoff - Fracture Vug Workflow","Standoff Porosity","Thermal Neutron Porosity","Total Porosity","Vug Porosity"]'
$ws.Range("B25").Value = $ws.Range("C24").Value2
$ws.Range("B26").Value = $ws.Range("C24").Value2
$ws.Range("B27").Value = $ws.Range("C24").Value2
$ws.Range("B28").Value = $ws.Range("C24").Value2
$ws.Range("B29").Value = $ws.Range("C24").Value2
$ws.Range("B35").Value = $ws.Range("C24").Value2
$ws.Range("B36").Value = $ws.Range("C24").Value2
$ws.Range("B24").Value = '["Core Permeability","Core Permeability Log10","Core Permeability Vertical","Core Permeability Under Stress (Array)","Average Permeability","Horizontal Permeability","Linear Permeability","Net Permeability","NMR Permeability","Permeability","Vertical Permeability","Permeability Pc Modeling"]'
$ws.Range("C25").Value = $ws.Range("B24").Value2
$ws.Range("C26").Value = $ws.Range("B24").Value2
$ws.Range("C27").Value = $ws.Range("B24").Value2
$ws.Range("C28").Value = $ws.Range("B24").Value2
$ws.Range("C29").Value = $ws.Range("B24").Value2

# --- Re-apply sort on the data range (fixes stale sort condition reference) ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A56"))
$sortObj.SetRange($ws.Range("A2:E56"))
$sortObj.Header = 2
$sortObj.Apply()

# --- Update selection / view ---
$ws.Range("C55").Select()
